$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Day 9 - Test Case Summary" block, mirroring the existing Day N blocks
# (header row + 3 data rows), appended after the Day 8 block (rows 48-51).

# Header row 54 (merged B54:C54) — same style/text pattern as the other headers.
$ws.Range("B54").Value = "Spint( 37) - Day 9 - Test Case Summary"
$ws.Range("B54:C54").Merge()

# Data rows 55-57
$ws.Range("B55").Value = "Total  testcase Written"
$ws.Range("C55").Value = 408

$ws.Range("B56").Value = "Total Execution"
$ws.Range("C56").Value = 639

$ws.Range("B57").Value = "Total Review"
$ws.Range("C57").Value = 511

# Copy formatting (styles + row heights) from the Day 8 block (rows 48-51)
# onto the new Day 9 block (rows 54-57).
$ws.Range("B48:C48").Copy() | Out-Null
$ws.Range("B54:C54").PasteSpecial(-4122) | Out-Null

$ws.Range("B49:C49").Copy() | Out-Null
$ws.Range("B55:C55").PasteSpecial(-4122) | Out-Null

$ws.Range("B50:C50").Copy() | Out-Null
$ws.Range("B56:C56").PasteSpecial(-4122) | Out-Null

$ws.Range("B51:C51").Copy() | Out-Null
$ws.Range("B57:C57").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Restore the text/values in case PasteSpecial(formats) disturbed anything.
$ws.Range("B54").Value = "Spint( 37) - Day 9 - Test Case Summary"
$ws.Range("B55").Value = "Total  testcase Written"
$ws.Range("C55").Value = 408
$ws.Range("B56").Value = "Total Execution"
$ws.Range("C56").Value = 639
$ws.Range("B57").Value = "Total Review"
$ws.Range("C57").Value = 511

# Updated "Total Execution" figure for the Day 8 block (C49: 369 -> 371)
$ws.Range("C49").Value = 371

# Row heights for the new block (match the other blocks' pattern:
# header/data/data rows = 18.75, the "Total testcase Written" row = 37.5).
$ws.Rows.Item(54).RowHeight = 18.75
$ws.Rows.Item(55).RowHeight = 37.5
$ws.Rows.Item(56).RowHeight = 18.75
$ws.Rows.Item(57).RowHeight = 18.75

# View state follow-on from the new rows being added (selection moves to
# the new last data cell, mirroring the original author's cursor position).
$ws.Range("C56").Select()
